$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '<<FORM(["Income":"Financial Aid"],AVERAGE)>>'

$ws.Range("E2").Select()
